$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has area rows in A2:B78 (A = index 1..77, B = area name)
# and value columns C:F populated through row 79 (row 79 only has C:F, no A/B yet).
# We need to insert a new "Chicago" entry at the top (row 2), shifting the
# existing A/B values down by one row each (C:F stay on their current row).
#
# Work from the bottom up so we don't clobber values before they are copied.
$lastRow = 79

for ($r = $lastRow; $r -ge 3; $r--) {
    $srcRow = $r - 1
    $aVal = $ws.Cells.Item($srcRow, 1).Value2
    $bVal = $ws.Cells.Item($srcRow, 2).Value2
    $ws.Cells.Item($r, 1).Value2 = $aVal
    $ws.Cells.Item($r, 2).Value2 = $bVal
}

# New top entry: Chicago (the overall city total), index 0
$ws.Cells.Item(2, 1).Value2 = 0
$ws.Cells.Item(2, 2).Value2 = "Chicago"

# Update the saved selection to B3, matching the edited file
$ws.Range("B3").Select()

$wb.Save()
